$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# K13 used to hold the old "ноутбук 3 фунта..." note; it is now blank
# (keeps its style, loses its value).
$ws.Range("K13").ClearContents()

# Bring row 14 into existence with the same look (styles + row height) as
# the rows above it, then fill in the new "IPHONE" knapsack item.
$srcRow = $ws.Range("D13:J13")
$dstRow = $ws.Range("D14:J14")
$srcRow.Copy()
$dstRow.PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Rows.Item(14).RowHeight = 50.1

$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 2000
$ws.Range("F14").Value = "IPHONE"
$ws.Range("G14").Value = "2000`nI"
$ws.Range("H14").Value = "3500`nI + г"
$ws.Range("I14").Value = "3500`nI + г"
$ws.Range("J14").Value = "4000`nI + н"

# Move the active selection to where the author left it.
[void]$ws.Range("H20").Select()
